$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Đơn sale chính" ---
$ws1 = $wb.Worksheets.Item("Đơn sale chính")

$headers = @("Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách","Tên dịch vụ","Đơn giá gốc","Sale phụ","Upsale","Đơn giá","Đã thanh toán","Tỉ lệ chiết khấu sale chính","Chiết khấu sale chính")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws1.Range("A2").Value = "HD-LUXURY"
$ws1.Range("B2").Value = 614
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "08-01-2024"
$ws1.Range("C2").Style = "Normal"
$ws1.Range("D2").Value = "CẦN THƠ"
$ws1.Range("E2").Value = "Trần Nguyễn Yến Linh"
$ws1.Range("F2").Value = "Khách cũ"
$ws1.Range("G2").Value = "Cắt mí"
$ws1.Range("H2").Value = 0
$ws1.Range("I2").Value = "Đỗ Thị Huyền Trân"
$ws1.Range("J2").Value = 6000000
$ws1.Range("K2").Value = 6000000
$ws1.Range("L2").Value = 6000000
$ws1.Range("M2").Value = 0.1
$ws1.Range("N2").Value = 360000

$ws1.Range("A3").Value = "Tổng"
$ws1.Range("B3").Value = 1
$ws1.Range("C3").Value = ""
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = ""
$ws1.Range("F3").Value = ""
$ws1.Range("G3").Value = ""
$ws1.Range("H3").Value = 0
$ws1.Range("I3").Value = ""
$ws1.Range("J3").Value = 6000000
$ws1.Range("K3").Value = 6000000
$ws1.Range("L3").Value = 6000000
$ws1.Range("M3").Value = 0
$ws1.Range("N3").Value = 360000

# --- Sheet 2: "Lương" ---
$ws2 = $wb.Worksheets.Item("Lương")

$ws2.Range("B2").Value = 1.5
$ws2.Range("B3").Value = 52500
$ws2.Range("B4").Value = 160714.2857142857
$ws2.Range("B5").Value = 360000
$ws2.Range("B14").Value = 160714.2857142857
$ws2.Range("B24").Value = 160714.2857142857
$ws2.Range("B32").Value = 573214.2857142857
$ws2.Range("B33").Value = 160714.2857142857
$ws2.Range("B34").Value = 160714.2857142857
$ws2.Range("A35").Value = "Tổng lương tại HỆ THỐNG"
$ws2.Range("B35").Value = 894642.857142857
